{"js": "const body = context.document.body;\n\nconst replacements = [\n    [\"2023-11-22 Wednesday\", \"2023-11-23 Thursday\"],\n    [\"54\u00d737=\", \"11\u00d752=\"],\n    [\"27\u00d740=\", \"57\u00d798=\"],\n    [\"48\u00d719=\", \"84\u00d754=\"],\n    [\"34\u00d741=\", \"24\u00d762=\"],\n    [\"39\u00d780=\", \"49\u00d722=\"],\n    [\"17\u00d788=\", \"50\u00d724=\"],\n    [\"35\u00d796=\", \"43\u00d714=\"],\n    [\"43\u00d779=\", \"21\u00d736=\"],\n    [\"93\u00d719=\", \"94\u00d797=\"],\n    [\"83\u00d745=\", \"42\u00d747=\"],\n    [\"98\u00d793=\", \"34\u00d781=\"],\n    [\"23\u00d767=\", \"26\u00d721=\"],\n    [\"31\u00d779=\", \"44\u00d725=\"],\n    [\"15\u00d774=\", \"61\u00d757=\"],\n    [\"65\u00d740=\", \"29\u00d724=\"],\n    [\"75\u00d722=\", \"55\u00d756=\"],\n    [\"34\u00d760=\", \"71\u00d735=\"],\n    [\"31\u00d751=\", \"47\u00d791=\"],\n    [\"37\u00d731=\", \"71\u00d782=\"],\n    [\"30\u00d712=\", \"65\u00d732=\"],\n    [\"12\u00d712=\", \"75\u00d712=\"],\n    [\"99\u00d794=\", \"55\u00d723=\"],\n    [\"46\u00d768=\", \"59\u00d764=\"],\n    [\"98\u00d774=\", \"95\u00d749=\"],\n    [\"94\u00d766=\", \"47\u00d764=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n    const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    await context.sync();\n    for (let i = 0; i < results.items.length; i++) {\n        results.items[i].insertText(newText, Word.InsertLocation.replace);\n    }\n    await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old = \"2023-11-22 Wednesday\"; New = \"2023-11-23 Thursday\"},\n    @{Old = \"54\u00d737=\"; New = \"11\u00d752=\"},\n    @{Old = \"27\u00d740=\"; New = \"57\u00d798=\"},\n    @{Old = \"48\u00d719=\"; New = \"84\u00d754=\"},\n    @{Old = \"34\u00d741=\"; New = \"24\u00d762=\"},\n    @{Old = \"39\u00d780=\"; New = \"49\u00d722=\"},\n    @{Old = \"17\u00d788=\"; New = \"50\u00d724=\"},\n    @{Old = \"35\u00d796=\"; New = \"43\u00d714=\"},\n    @{Old = \"43\u00d779=\"; New = \"21\u00d736=\"},\n    @{Old = \"93\u00d719=\"; New = \"94\u00d797=\"},\n    @{Old = \"83\u00d745=\"; New = \"42\u00d747=\"},\n    @{Old = \"98\u00d793=\"; New = \"34\u00d781=\"},\n    @{Old = \"23\u00d767=\"; New = \"26\u00d721=\"},\n    @{Old = \"31\u00d779=\"; New = \"44\u00d725=\"},\n    @{Old = \"15\u00d774=\"; New = \"61\u00d757=\"},\n    @{Old = \"65\u00d740=\"; New = \"29\u00d724=\"},\n    @{Old = \"75\u00d722=\"; New = \"55\u00d756=\"},\n    @{Old = \"34\u00d760=\"; New = \"71\u00d735=\"},\n    @{Old = \"31\u00d751=\"; New = \"47\u00d791=\"},\n    @{Old = \"37\u00d731=\"; New = \"71\u00d782=\"},\n    @{Old = \"30\u00d712=\"; New = \"65\u00d732=\"},\n    @{Old = \"12\u00d712=\"; New = \"75\u00d712=\"},\n    @{Old = \"99\u00d794=\"; New = \"55\u00d723=\"},\n    @{Old = \"46\u00d768=\"; New = \"59\u00d764=\"},\n    @{Old = \"98\u00d774=\"; New = \"95\u00d749=\"},\n    @{Old = \"94\u00d766=\"; New = \"47\u00d764=\"}\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Execute([ref]$pair.Old, $false, $true, $false, $false, $false, $true, 1, $false, [ref]$pair.New, 2)\n}\n"}
